$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.155.26"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.853.12"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.97%  "
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "310.38"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "0.4776"
$ws.Range("D8").Value = "0.3697"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "0.07279"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "0.9333"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "0.07809"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "1.851.88"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "5.389"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "89.60"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "1.020"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "0.000008708"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D20").Value = "27.146.09"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "14.63"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "5.076"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "10.66"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "1.940"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "153.05"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").Value = "1.986"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").Value = "115.03"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "4.930"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "0.08875"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("D32").Value = "1.181"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "4.535"
$ws.Range("D34").Value = "0.7382"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").Value = "0.01991"
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("D39").Value = "2.983"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").Value = "0.5295"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").Value = "7.052"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").Value = "0.1527"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").Value = "8.312"
$ws.Range("D44").Value = "10.56"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "0.4749"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").Value = "102.15"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "1.626"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "65.90"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").Value = "0.06058"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "0.8942"
$ws.Range("E51").Value = "  +0.55%  "
